$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 - Shuttle
$ws.Range("A15").Value = "Shuttle"
$ws.Range("B15").Value = 57999
$ws.Range("C15").Value = 9
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = 45586
$ws.Range("F15").Value = 49
$ws.Range("G15").Value = 171
$ws.Range("H15").Value = 8903
$ws.Range("I15").Value = 3267
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = 13

# Row 16 - Penbased
$ws.Range("A16").Value = "Penbased"
$ws.Range("B16").Value = 10992
$ws.Range("C16").Value = 16
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 1143
$ws.Range("F16").Value = 1143
$ws.Range("G16").Value = 1144
$ws.Range("H16").Value = 1055
$ws.Range("I16").Value = 1144
$ws.Range("J16").Value = 1055
$ws.Range("K16").Value = 1056
$ws.Range("L16").Value = 1142
$ws.Range("M16").Value = 1055
$ws.Range("N16").Value = 1055

# Row 17 - Magic
$ws.Range("A17").Value = "Magic"
$ws.Range("B17").Value = 19020
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 12332
$ws.Range("F17").Value = 6688

# Column B on the new rows carries the same left-border style as B14
$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the last-used selection recorded in the saved file
$ws.Range("J23").Select() | Out-Null
